$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 26.13880918633845
    "C2" = 6.494139229403823
    "D2" = 13.3469647552972
    "E2" = 13.25214728482555
    "G2" = 72.70950745627054
    "H2" = 24.92428454187738
    "J2" = 8.11127083859874
    "L2" = 13.55151225906293
    "B3" = 25.91458437953561
    "C3" = 6.1684309188031
    "D3" = 13.353890004925
    "E3" = 13.2768983479744
    "G3" = 72.27631281579214
    "H3" = 24.90078310381657
    "J3" = 8.114537807367885
    "L3" = 13.5615865156254
    "B4" = 25.78388185009584
    "C4" = 5.961823492362269
    "D4" = 13.36046839528703
    "E4" = 13.2934550553058
    "G4" = 72.02621939206459
    "H4" = 24.89090967871217
    "J4" = 8.116665879284195
    "L4" = 13.57024054995178
    "B5" = 25.73242891012261
    "C5" = 5.876102549606684
    "D5" = 13.36373360340621
    "E5" = 13.30054404691454
    "G5" = 71.92835436858064
    "H5" = 24.88803053075184
    "J5" = 8.117563872930416
    "L5" = 13.57438775753321
    "B6" = 25.72399604050945
    "C6" = 5.861780902048117
    "D6" = 13.3643110706968
    "E6" = 13.30174182895347
    "G6" = 71.91235002309153
    "H6" = 24.88762151105479
    "J6" = 8.117714845738435
    "L6" = 13.5751138762062
    "B7" = 25.78318054024014
    "C7" = 5.960673413059384
    "D7" = 13.36051006535423
    "E7" = 13.29354927498353
    "G7" = 72.02488308512918
    "H7" = 24.89086621865951
    "J7" = 8.116677865185441
    "L7" = 13.57029396803922
    "B8" = 26.06008411867064
    "C8" = 6.383285724526318
    "D8" = 13.34886951693235
    "E8" = 13.26039938899866
    "G8" = 72.5568787635147
    "H8" = 24.91523486849
    "J8" = 8.112371983784072
    "L8" = 13.55447354534415
    "B9" = 26.65578070353837
    "C9" = 7.154431876375974
    "D9" = 13.34451689747253
    "E9" = 13.20617268167571
    "G9" = 73.72349465669163
    "H9" = 24.99921448655355
    "J9" = 8.104893949193247
    "L9" = 13.54303375754969
    "B10" = 27.12182877788326
    "C10" = 7.680342478930971
    "D10" = 13.35259616507987
    "E10" = 13.1728956479996
    "G10" = 74.65177573189811
    "H10" = 25.08299467161973
    "J10" = 8.099983986123915
    "L10" = 13.54655790223001
    "B11" = 27.3391704826802
    "C11" = 7.909897504619316
    "D11" = 13.35871965483867
    "E11" = 13.15918046869289
    "G11" = 75.08854530258458
    "H11" = 25.12589184531476
    "J11" = 8.097876187810597
    "L11" = 13.55074552117401
    "B12" = 27.42216657561005
    "C12" = 7.99687350835382
    "D12" = 13.36139003914584
    "E12" = 13.1541913398096
    "G12" = 75.25593120345235
    "H12" = 25.14282163300145
    "J12" = 8.097096031848443
    "L12" = 13.55270202258324
    "B13" = 27.40426221655731
    "C13" = 7.977031430672795
    "D13" = 13.3607992988544
    "E13" = 13.1552567434223
    "G13" = 75.2197946409735
    "H13" = 25.13914505813807
    "J13" = 8.097263251888403
    "L13" = 13.55226418272032
    "B14" = 27.34598514610227
    "C14" = 7.916959161501602
    "D14" = 13.35893230691989
    "E14" = 13.15876591175144
    "G14" = 75.10227689145481
    "H14" = 25.12727094669696
    "J14" = 8.09781164309396
    "L14" = 13.55089905941203
    "B15" = 27.31037687743699
    "C15" = 7.879972852809016
    "D15" = 13.35783448460697
    "E15" = 13.16094201000074
    "G15" = 75.03055029864994
    "H15" = 25.12008690258271
    "J15" = 8.098149893597146
    "L15" = 13.55011113314097
    "B16" = 27.10772605543846
    "C16" = 7.665140526279981
    "D16" = 13.35224520689505
    "E16" = 13.17382058964755
    "G16" = 74.62351554886425
    "H16" = 25.08028732967389
    "J16" = 8.100124261004826
    "L16" = 13.546336165089
    "B17" = 26.98472049681468
    "C17" = 7.530823473677331
    "D17" = 13.34944312001444
    "E17" = 13.18208552032937
    "G17" = 74.37745947416661
    "H17" = 25.05709560052975
    "J17" = 8.101367638055693
    "L17" = 13.54468171110355
    "B18" = 26.91447749467547
    "C18" = 7.452659009520706
    "D18" = 13.34806191189453
    "E18" = 13.18697320649488
    "G18" = 74.23730424616963
    "H18" = 25.04420693155105
    "J18" = 8.102094637335588
    "L18" = 13.54397344220581
    "B19" = 26.89078364633079
    "C19" = 7.426039681359192
    "D19" = 13.34763385506536
    "E19" = 13.18865109552836
    "G19" = 74.19008810036425
    "H19" = 25.03992050583933
    "J19" = 8.102342822582125
    "L19" = 13.54377545070503
    "B20" = 26.99776275344604
    "C20" = 7.545216261947625
    "D20" = 13.349717556266
    "E20" = 13.18119184500882
    "G20" = 74.403511461414
    "H20" = 25.05951778033407
    "J20" = 8.101234053342523
    "L20" = 13.54483265279631
    "B21" = 27.36308430244624
    "C21" = 7.934643591182017
    "D21" = 13.35947115167568
    "E21" = 13.15772963413734
    "G21" = 75.13674143678918
    "H21" = 25.13074007907844
    "J21" = 8.097650078709338
    "L21" = 13.55128997561268
    "B22" = 27.60585081656965
    "C22" = 8.210992115702572
    "D22" = 13.36789442541768
    "E22" = 13.14358774042579
    "G22" = 75.62750715656904
    "H22" = 25.18128223329221
    "J22" = 8.095412758350012
    "L22" = 13.55767060951065
    "B23" = 27.47593950736275
    "C23" = 8.06033040914634
    "D23" = 13.3632115236098
    "E23" = 13.15102648966904
    "G23" = 75.36455052222004
    "H23" = 25.15394248820695
    "J23" = 8.096597270360133
    "L23" = 13.55406781893658
    "B24" = 26.9918648655948
    "C24" = 7.538712217665863
    "D24" = 13.34959276798494
    "E24" = 13.18159545179603
    "G24" = 74.39172928406535
    "H24" = 25.05842132793992
    "J24" = 8.101294409094191
    "L24" = 13.54476365541569
    "B25" = 26.48936371206138
    "C25" = 6.952616082744317
    "D25" = 13.34371429069615
    "E25" = 13.21968923519066
    "G25" = 73.39507930445308
    "H25" = 24.97261650779716
    "J25" = 8.106814042713074
    "L25" = 13.54403218868781
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
